$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-14 04:43:42"
$wsZhCn.Range("H2").Value = "2016-03-14 04:43:59"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-14 04:43:45"
$wsDeDe.Range("H2").Value = "2016-03-14 04:44:05"
